$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'91.144.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.29%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.144.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.61%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.14%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'238.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +9.30%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'635.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.75%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.08"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.53%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.368"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.09%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.00%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'3.143.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +0.71%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.727"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -5.56%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.82%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'36.76"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +4.94%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -1.18%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.52%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'90.916.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.28%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.721.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.61%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.172.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.28%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -4.24%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'Chainlink"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'14.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.18%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "'PEPE"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'0.0000213"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.97%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'447.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.85%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +9.84%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'9.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.45%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'6.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -5.68%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'90.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.05%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'12.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.81%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'3.307.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.70%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.01%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'9.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +3.74%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -4.02%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.975"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +8.66%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'27.13"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +11.64%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.198"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +23.94%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.68%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'515.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.89%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +2.74%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'7.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.15%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +3.84%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'1.31"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +0.21%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.421"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +5.02%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.35%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.0855"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -3.23%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.00%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +47.30%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +0.10%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'150.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.89%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.696"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +9.81%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'ImmutableX"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'1.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +3.97%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'OKB"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'45.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +3.36%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'Filecoin"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'4.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +5.36%  "
$ws.Range("E51").Style = "Normal"
